# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.590.08'
$ws.Range("E2").Value = '  +1.91%  '
$ws.Range("D3").Value = '2.161.66'
$ws.Range("E3").Value = '  +0.87%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '227.02'
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("E6").Value = '  +0.76%  '
$ws.Range("D7").Value = '62.68'
$ws.Range("E7").Value = '  +0.74%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = '0.0847'
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("E11").Value = '  +0.54%  '
$ws.Range("D12").Value = '15.86'
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("D13").Value = '2.482.78'
$ws.Range("E13").Value = '  +1.10%  '
$ws.Range("D14").Value = '21.72'
$ws.Range("E14").Value = '  -1.99%  '
$ws.Range("E15").Value = '  -0.66%  '
$ws.Range("D16").Value = '5.46'
$ws.Range("E16").Value = '  -0.98%  '
$ws.Range("D17").Value = '2.163.07'
$ws.Range("E17").Value = '  +1.26%  '
$ws.Range("D18").Value = '39.576.48'
$ws.Range("E18").Value = '  +1.82%  '
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").Value = '6.03'
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("D21").Value = '0.0₃0880'
$ws.Range("E21").Value = '  +3.99%  '
$ws.Range("D22").Value = '227.74'
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("E25").Value = '  -2.08%  '
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("D27").Value = '9.43'
$ws.Range("E27").Value = '  -1.22%  '
$ws.Range("E28").Value = '  +0.62%  '
$ws.Range("E29").Value = '  +0.78%  '
$ws.Range("D30").Value = '19.62'
$ws.Range("E30").Value = '  +0.34%  '
$ws.Range("E31").Value = '  +4.83%  '
$ws.Range("D33").Value = '4.55'
$ws.Range("E33").Value = '  -0.82%  '
$ws.Range("E34").Value = '  -2.32%  '
$ws.Range("D35").Value = '6.96'
$ws.Range("E35").Value = '  -2.70%  '
$ws.Range("E36").Value = '  +0.29%  '
$ws.Range("D37").Value = '3.83'
$ws.Range("E37").Value = '  +7.96%  '
$ws.Range("E38").Value = '  -0.86%  '
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  +0.07%  '
$ws.Range("E40").Value = '  +19.62%  '
$ws.Range("D41").Value = '102.42'
$ws.Range("E41").Value = '  -0.22%  '
$ws.Range("E42").Value = '  -1.14%  '
$ws.Range("D43").Value = '17.64'
$ws.Range("E43").Value = '  -2.86%  '
$ws.Range("D44").Value = '1.512.81'
$ws.Range("E44").Value = '  -1.40%  '
$ws.Range("D45").Value = '1.21'
$ws.Range("E45").Value = '  +1.97%  '
$ws.Range("E46").Value = '  +0.59%  '
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("E48").Value = '  +0.13%  '
$ws.Range("E49").Value = '  -0.82%  '
$ws.Range("E50").Value = '  +28.19%  '
$ws.Range("E51").Value = '  +0.35%  '
